$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "23.469.33"
Set-TextCell "E2" "  +1.21%  "
Set-TextCell "D3" "1.638.56"
Set-TextCell "E3" "  +2.25%  "
Set-TextCell "D4" "1.002"
Set-TextCell "E4" "  +0.11%  "
Set-TextCell "D5" "307.19"
Set-TextCell "E5" "  +1.40%  "
Set-TextCell "D6" "1.002"
Set-TextCell "E6" "  +0.08%  "
Set-TextCell "D7" "0.3771"
Set-TextCell "E7" "  -0.33%  "
Set-TextCell "D8" "52.19"
Set-TextCell "E8" "  -0.22%  "
Set-TextCell "D9" "0.3646"
Set-TextCell "E9" "  +0.85%  "
Set-TextCell "D10" "1.267"
Set-TextCell "E10" "  -0.01%  "
Set-TextCell "D11" "0.08170"
Set-TextCell "E11" "  +0.46%  "
Set-TextCell "D12" "1.002"
Set-TextCell "E12" "  +0.13%  "
Set-TextCell "D13" "23.01"
Set-TextCell "E13" "  +1.56%  "
Set-TextCell "D14" "6.642"
Set-TextCell "E14" "  +0.83%  "
Set-TextCell "D15" "0.00001279"
Set-TextCell "E15" "  +2.52%  "
Set-TextCell "D16" "7.385"
Set-TextCell "E16" "  -0.13%  "
Set-TextCell "D17" "1.644.03"
Set-TextCell "E17" "  +2.75%  "
Set-TextCell "D18" "94.87"
Set-TextCell "E18" "  +1.08%  "
Set-TextCell "D19" "0.06964"
Set-TextCell "E19" "  +1.20%  "
Set-TextCell "D20" "18.21"
Set-TextCell "E20" "  +0.78%  "
Set-TextCell "D21" "6.552"
Set-TextCell "E21" "  -0.02%  "
Set-TextCell "E22" "  -0.01%  "
Set-TextCell "D23" "23.465.33"
Set-TextCell "E23" "  +1.16%  "
Set-TextCell "D24" "12.81"
Set-TextCell "E24" "  -1.02%  "
Set-TextCell "D26" "2.423"
Set-TextCell "E26" "  +1.54%  "
Set-TextCell "D27" "21.28"
Set-TextCell "E27" "  +0.31%  "
Set-TextCell "D28" "151.38"
Set-TextCell "E28" "  +1.61%  "
Set-TextCell "D29" "5.362"
Set-TextCell "E29" "  +2.05%  "
Set-TextCell "D30" "135.64"
Set-TextCell "E30" "  +1.35%  "
Set-TextCell "D31" "2.337"
Set-TextCell "E31" "  -1.55%  "
Set-TextCell "D32" "1.824.36"
Set-TextCell "E32" "  +2.63%  "
Set-TextCell "D33" "6.789"
Set-TextCell "E33" "  -0.64%  "
Set-TextCell "D34" "0.9653"
Set-TextCell "E34" "  -0.81%  "
Set-TextCell "D35" "0.02829"
Set-TextCell "E35" "  +4.08%  "
Set-TextCell "D36" "10.32"
Set-TextCell "E36" "  -0.04%  "
Set-TextCell "D37" "0.07359"
Set-TextCell "E37" "  -2.16%  "
Set-TextCell "D38" "0.2537"
Set-TextCell "E38" "  +1.06%  "
Set-TextCell "D39" "6.173"
Set-TextCell "E39" "  +0.78%  "
Set-TextCell "D40" "0.08859"
Set-TextCell "E40" "  +0.53%  "
Set-TextCell "E41" "  +1.75%  "
Set-TextCell "D42" "0.7110"
Set-TextCell "E42" "  +0.12%  "
Set-TextCell "E43" "  -0.05%  "
Set-TextCell "D44" "16.28"
Set-TextCell "E44" "  +4.82%  "
Set-TextCell "D45" "0.6558"
Set-TextCell "E45" "  +0.34%  "
Set-TextCell "D46" "2.341"
Set-TextCell "E46" "  +1.16%  "
Set-TextCell "D47" "1.001"
Set-TextCell "E47" "  +0.07%  "
Set-TextCell "D48" "4.028"
Set-TextCell "E48" "  +0.43%  "
Set-TextCell "D49" "0.07971"
Set-TextCell "E49" "  +0.25%  "
Set-TextCell "D50" "129.45"
Set-TextCell "E50" "  -2.03%  "
Set-TextCell "D51" "1.208"
Set-TextCell "E51" "  +0.39%  "
